# Update "想去人数" (want-to-go count) figures in column F
# for the 展览 (Exhibitions) and 全部类型 (All Types) sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 8208
$wsExhibit.Range("F5").Value = 5981
$wsExhibit.Range("F6").Value = 512
$wsExhibit.Range("F10").Value = 304
$wsExhibit.Range("F11").Value = 797
$wsExhibit.Range("F12").Value = 77

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 8208
$wsAll.Range("F5").Value = 5981
$wsAll.Range("F6").Value = 512
$wsAll.Range("F10").Value = 304
$wsAll.Range("F15").Value = 797
$wsAll.Range("F16").Value = 77
